# Apply an AutoFilter on the "IO Type" column (3rd column, colId=2) of
# Table1 in Sheet1, keeping only rows whose IO Type is "DO". This filters
# out (hides) every row that doesn't match, same as a user selecting the
# "DO" checkbox in the column's filter dropdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$lo = $ws.ListObjects.Item("Table1")

# Field 3 == "IO Type" column within the table (1-based).
# Passing the criteria as an array produces a standard value-list filter
# (<filters><filter val="DO"/></filters>) rather than a custom filter.
$lo.Range.AutoFilter(3, @("DO"))
